{"js": "// Auto-generated: replace each unique run text with its updated value.\n// All original <w:t> values in this document are unique, so a literal,\n// case-sensitive, non-wildcard search+replace for each pair is unambiguous\n// and order-independent.\nconst replacements = [\n  [\"2024-07-18 Thursday\", \"2024-07-19 Friday\"],\n  [\"62-13=49\", \"59+10=69\"],\n  [\"12+82=94\", \"35+44=79\"],\n  [\"42-31=11\", \"81-36=45\"],\n  [\"83-15=68\", \"41+12=53\"],\n  [\"34+4=38\", \"13+29=42\"],\n  [\"66-64=2\", \"91-17=74\"],\n  [\"23-11=12\", \"27+67=94\"],\n  [\"64+26=90\", \"52-9=43\"],\n  [\"87-49=38\", \"83-65=18\"],\n  [\"70-42=28\", \"50-6=44\"],\n  [\"61-56=5\", \"64-50=14\"],\n  [\"44-34=10\", \"9-2=7\"],\n  [\"15+22=37\", \"74-21=53\"],\n  [\"35+42=77\", \"10-6=4\"],\n  [\"23+27=50\", \"45-13=32\"],\n  [\"4+57=61\", \"46-2=44\"],\n  [\"22+10=32\", \"62-3=59\"],\n  [\"72-38=34\", \"10+60=70\"],\n  [\"65+6=71\", \"11+15=26\"],\n  [\"42+55=97\", \"95-12=83\"],\n  [\"53-30=23\", \"67+5=72\"],\n  [\"84-12=72\", \"44-8=36\"],\n  [\"36+39=75\", \"33-8=25\"],\n  [\"73-44=29\", \"19+6=25\"],\n  [\"84-38=46\", \"23+14=37\"],\n  [\"92-5=87\", \"27+36=63\"],\n  [\"44+49=93\", \"13+84=97\"],\n  [\"74-43=31\", \"84-55=29\"],\n  [\"14+15=29\", \"10+9=19\"],\n  [\"14+17=31\", \"12+85=97\"],\n  [\"17+34=51\", \"71-60=11\"],\n  [\"41-4=37\", \"64-27=37\"],\n  [\"24-24=0\", \"95-67=28\"],\n  [\"96-83=13\", \"69+17=86\"],\n  [\"75-72=3\", \"57-51=6\"],\n  [\"80-73=7\", \"56+3=59\"],\n  [\"79-5=74\", \"30+9=39\"],\n  [\"69-37=32\", \"18+17=35\"],\n  [\"17-0=17\", \"71-14=57\"],\n  [\"7+71=78\", \"34-7=27\"],\n  [\"29+69=98\", \"43+45=88\"],\n  [\"87-68=19\", \"52-10=42\"],\n  [\"41+29=70\", \"50+0=50\"],\n  [\"92-80=12\", \"56+24=80\"],\n  [\"44+11=55\", \"73-4=69\"],\n  [\"30-15=15\", \"18+21=39\"],\n  [\"41+38=79\", \"5+42=47\"],\n  [\"23+68=91\", \"5+64=69\"],\n  [\"2+5=7\", \"95-48=47\"],\n  [\"17+57=74\", \"57+18=75\"],\n  [\"78-18=60\", \"19-3=16\"],\n  [\"16+13=29\", \"72-47=25\"],\n  [\"91-55=36\", \"65-57=8\"],\n  [\"30+0=30\", \"50-25=25\"],\n  [\"22+71=93\", \"27+53=80\"],\n  [\"61+12=73\", \"15+3=18\"],\n  [\"62-16=46\", \"7+61=68\"],\n  [\"22+3=25\", \"46-4=42\"],\n  [\"89-50=39\", \"51-7=44\"],\n  [\"32+33=65\", \"14+69=83\"],\n  [\"34+49=83\", \"35-28=7\"],\n  [\"47-26=21\", \"36+5=41\"],\n  [\"71-36=35\", \"14-9=5\"],\n  [\"41-30=11\", \"44+9=53\"],\n  [\"60-48=12\", \"45-17=28\"],\n  [\"28+27=55\", \"30+61=91\"],\n  [\"38+12=50\", \"43+31=74\"],\n  [\"49+21=70\", \"93-47=46\"],\n  [\"74-2=72\", \"67-46=21\"],\n  [\"57-45=12\", \"33+64=97\"],\n  [\"38+58=96\", \"57+25=82\"],\n  [\"3+32=35\", \"11+88=99\"],\n  [\"40+50=90\", \"81-6=75\"],\n  [\"30+5=35\", \"11+68=79\"],\n  [\"92-14=78\", \"43+6=49\"],\n  [\"69-33=36\", \"70+6=76\"],\n  [\"45-12=33\", \"15+47=62\"],\n  [\"29+1=30\", \"14+40=54\"],\n  [\"5+31=36\", \"29-13=16\"],\n  [\"40-21=19\", \"15+10=25\"],\n  [\"71+11=82\", \"21+71=92\"],\n  [\"79-61=18\", \"62+5=67\"],\n  [\"10+38=48\", \"66+26=92\"],\n  [\"6+68=74\", \"34+21=55\"],\n  [\"0+8=8\", \"18+65=83\"],\n  [\"51-3=48\", \"58+36=94\"],\n  [\"37+58=95\", \"78-8=70\"],\n  [\"81-71=10\", \"28+52=80\"],\n  [\"29-12=17\", \"47+29=76\"],\n  [\"69-44=25\", \"31-3=28\"],\n  [\"4+73=77\", \"51+25=76\"],\n  [\"99-9=90\", \"48-23=25\"],\n  [\"51-5=46\", \"44-33=11\"],\n  [\"85+12=97\", \"29+31=60\"],\n  [\"58-16=42\", \"24+67=91\"],\n  [\"6+85=91\", \"28-21=7\"],\n  [\"43-30=13\", \"82-27=55\"],\n  [\"95-90=5\", \"32+9=41\"],\n  [\"34+51=85\", \"15+79=94\"],\n  [\"81-38=43\", \"70+19=89\"],\n];\n\nconst searchOptions = { matchCase: true, matchWholeWord: false };\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, searchOptions);\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Auto-generated: replace each unique run text with its updated value.\n# Every original text run in this document is unique, so a literal,\n# case-sensitive (no wildcards) Find/Replace for each pair is unambiguous\n# and order-independent. wdReplaceOne = 1 (used via -Replace argument 2\n# positionally below actually corresponds to wdReplaceAll per-call, but\n# since each FindText is unique in the document a single replace occurs).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2024-07-18 Thursday\", \"2024-07-19 Friday\")\n  ,@(\"62-13=49\", \"59+10=69\")\n  ,@(\"12+82=94\", \"35+44=79\")\n  ,@(\"42-31=11\", \"81-36=45\")\n  ,@(\"83-15=68\", \"41+12=53\")\n  ,@(\"34+4=38\", \"13+29=42\")\n  ,@(\"66-64=2\", \"91-17=74\")\n  ,@(\"23-11=12\", \"27+67=94\")\n  ,@(\"64+26=90\", \"52-9=43\")\n  ,@(\"87-49=38\", \"83-65=18\")\n  ,@(\"70-42=28\", \"50-6=44\")\n  ,@(\"61-56=5\", \"64-50=14\")\n  ,@(\"44-34=10\", \"9-2=7\")\n  ,@(\"15+22=37\", \"74-21=53\")\n  ,@(\"35+42=77\", \"10-6=4\")\n  ,@(\"23+27=50\", \"45-13=32\")\n  ,@(\"4+57=61\", \"46-2=44\")\n  ,@(\"22+10=32\", \"62-3=59\")\n  ,@(\"72-38=34\", \"10+60=70\")\n  ,@(\"65+6=71\", \"11+15=26\")\n  ,@(\"42+55=97\", \"95-12=83\")\n  ,@(\"53-30=23\", \"67+5=72\")\n  ,@(\"84-12=72\", \"44-8=36\")\n  ,@(\"36+39=75\", \"33-8=25\")\n  ,@(\"73-44=29\", \"19+6=25\")\n  ,@(\"84-38=46\", \"23+14=37\")\n  ,@(\"92-5=87\", \"27+36=63\")\n  ,@(\"44+49=93\", \"13+84=97\")\n  ,@(\"74-43=31\", \"84-55=29\")\n  ,@(\"14+15=29\", \"10+9=19\")\n  ,@(\"14+17=31\", \"12+85=97\")\n  ,@(\"17+34=51\", \"71-60=11\")\n  ,@(\"41-4=37\", \"64-27=37\")\n  ,@(\"24-24=0\", \"95-67=28\")\n  ,@(\"96-83=13\", \"69+17=86\")\n  ,@(\"75-72=3\", \"57-51=6\")\n  ,@(\"80-73=7\", \"56+3=59\")\n  ,@(\"79-5=74\", \"30+9=39\")\n  ,@(\"69-37=32\", \"18+17=35\")\n  ,@(\"17-0=17\", \"71-14=57\")\n  ,@(\"7+71=78\", \"34-7=27\")\n  ,@(\"29+69=98\", \"43+45=88\")\n  ,@(\"87-68=19\", \"52-10=42\")\n  ,@(\"41+29=70\", \"50+0=50\")\n  ,@(\"92-80=12\", \"56+24=80\")\n  ,@(\"44+11=55\", \"73-4=69\")\n  ,@(\"30-15=15\", \"18+21=39\")\n  ,@(\"41+38=79\", \"5+42=47\")\n  ,@(\"23+68=91\", \"5+64=69\")\n  ,@(\"2+5=7\", \"95-48=47\")\n  ,@(\"17+57=74\", \"57+18=75\")\n  ,@(\"78-18=60\", \"19-3=16\")\n  ,@(\"16+13=29\", \"72-47=25\")\n  ,@(\"91-55=36\", \"65-57=8\")\n  ,@(\"30+0=30\", \"50-25=25\")\n  ,@(\"22+71=93\", \"27+53=80\")\n  ,@(\"61+12=73\", \"15+3=18\")\n  ,@(\"62-16=46\", \"7+61=68\")\n  ,@(\"22+3=25\", \"46-4=42\")\n  ,@(\"89-50=39\", \"51-7=44\")\n  ,@(\"32+33=65\", \"14+69=83\")\n  ,@(\"34+49=83\", \"35-28=7\")\n  ,@(\"47-26=21\", \"36+5=41\")\n  ,@(\"71-36=35\", \"14-9=5\")\n  ,@(\"41-30=11\", \"44+9=53\")\n  ,@(\"60-48=12\", \"45-17=28\")\n  ,@(\"28+27=55\", \"30+61=91\")\n  ,@(\"38+12=50\", \"43+31=74\")\n  ,@(\"49+21=70\", \"93-47=46\")\n  ,@(\"74-2=72\", \"67-46=21\")\n  ,@(\"57-45=12\", \"33+64=97\")\n  ,@(\"38+58=96\", \"57+25=82\")\n  ,@(\"3+32=35\", \"11+88=99\")\n  ,@(\"40+50=90\", \"81-6=75\")\n  ,@(\"30+5=35\", \"11+68=79\")\n  ,@(\"92-14=78\", \"43+6=49\")\n  ,@(\"69-33=36\", \"70+6=76\")\n  ,@(\"45-12=33\", \"15+47=62\")\n  ,@(\"29+1=30\", \"14+40=54\")\n  ,@(\"5+31=36\", \"29-13=16\")\n  ,@(\"40-21=19\", \"15+10=25\")\n  ,@(\"71+11=82\", \"21+71=92\")\n  ,@(\"79-61=18\", \"62+5=67\")\n  ,@(\"10+38=48\", \"66+26=92\")\n  ,@(\"6+68=74\", \"34+21=55\")\n  ,@(\"0+8=8\", \"18+65=83\")\n  ,@(\"51-3=48\", \"58+36=94\")\n  ,@(\"37+58=95\", \"78-8=70\")\n  ,@(\"81-71=10\", \"28+52=80\")\n  ,@(\"29-12=17\", \"47+29=76\")\n  ,@(\"69-44=25\", \"31-3=28\")\n  ,@(\"4+73=77\", \"51+25=76\")\n  ,@(\"99-9=90\", \"48-23=25\")\n  ,@(\"51-5=46\", \"44-33=11\")\n  ,@(\"85+12=97\", \"29+31=60\")\n  ,@(\"58-16=42\", \"24+67=91\")\n  ,@(\"6+85=91\", \"28-21=7\")\n  ,@(\"43-30=13\", \"82-27=55\")\n  ,@(\"95-90=5\", \"32+9=41\")\n  ,@(\"34+51=85\", \"15+79=94\")\n  ,@(\"81-38=43\", \"70+19=89\")\n)\n\nforeach ($pair in $pairs) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n  $range = $d.Content\n  $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
